$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.875.13"
$ws.Range("E2").Value = "  -3.61%  "

$ws.Range("D3").Value = "3.528.21"
$ws.Range("E3").Value = "  -3.96%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'604.59"
$ws.Range("E5").Value = "  -5.76%  "

$ws.Range("D6").Value = "'154.19"
$ws.Range("E6").Value = "  -3.24%  "

$ws.Range("D7").Value = "3.525.73"
$ws.Range("E7").Value = "  -3.98%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "'0.484"
$ws.Range("E9").Value = "  -2.66%  "

$ws.Range("D10").Value = "'0.141"
$ws.Range("E10").Value = "  -2.29%  "

$ws.Range("E11").Value = "  -3.57%  "

$ws.Range("E12").Value = "  -3.83%  "

$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = "  -4.49%  "

$ws.Range("D14").Value = "4.131.17"
$ws.Range("E14").Value = "  -3.80%  "

$ws.Range("D15").Value = "'32.02"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("D16").Value = "3.542.29"
$ws.Range("E16").Value = "  -4.02%  "

$ws.Range("D17").Value = "66.953.00"
$ws.Range("E17").Value = "  -3.47%  "

$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("E19").Value = "  -1.82%  "

$ws.Range("D20").Value = "'15.44"
$ws.Range("E20").Value = "  -3.16%  "

$ws.Range("D21").Value = "'451.02"
$ws.Range("E21").Value = "  -3.25%  "

$ws.Range("D22").Value = "'9.35"
$ws.Range("E22").Value = "  -5.59%  "

$ws.Range("D23").Value = "'0.637"
$ws.Range("E23").Value = "  -1.60%  "

$ws.Range("D24").Value = "'78.91"
$ws.Range("E24").Value = "  -0.43%  "

$ws.Range("D25").Value = "3.674.33"
$ws.Range("E25").Value = "  -3.82%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "'0.0000123"
$ws.Range("E27").Value = "  -2.11%  "

$ws.Range("D28").Value = "'10.20"
$ws.Range("E28").Value = "  -6.31%  "

$ws.Range("D29").Value = "'8.30"
$ws.Range("E29").Value = "  -8.00%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.54"
$ws.Range("E30").Value = "  -3.08%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.68"
$ws.Range("E31").Value = "  -1.53%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").Value = "'25.90"
$ws.Range("E33").Value = "  -3.61%  "

$ws.Range("D34").Value = "'1.89"
$ws.Range("E34").Value = "  -5.38%  "

$ws.Range("E35").Value = "  -3.96%  "

$ws.Range("E36").Value = "  -4.34%  "

$ws.Range("D37").Value = "3.524.11"
$ws.Range("E37").Value = "  -3.90%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("D41").Value = "'176.46"
$ws.Range("E41").Value = "  -0.47%  "

$ws.Range("E42").Value = "  -1.80%  "

$ws.Range("D43").Value = "'5.59"
$ws.Range("E43").Value = "  -4.90%  "

$ws.Range("D44").Value = "'0.0875"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("D45").Value = "'0.891"
$ws.Range("E45").Value = "  -3.73%  "

$ws.Range("D46").Value = "'45.76"
$ws.Range("E46").Value = "  -2.13%  "

$ws.Range("D47").Value = "'28.25"
$ws.Range("E47").Value = "  +2.81%  "

$ws.Range("D48").Value = "'2.69"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("E49").Value = "  -1.27%  "

$ws.Range("E50").Value = "  -3.25%  "

$ws.Range("D51").Value = "'7.64"
$ws.Range("E51").Value = "  -2.37%  "
